$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the "Details of implementation" rectangle shape (id 7) robustly
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -like "Details of implementation*") {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(1)
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- Edit 1: "Dependency injection" -> "Ninject" (paragraph 4) ---
$para4 = $tr.Paragraphs(4)
$para4.Text = "Ninject"

# --- Edit 2: paragraph 7 ---
# Before: "Async" + "/await " + "при работе " + "CRUD " + "операций"
# After : "Асинхронные " + "CRUD " + "операци" + "и"
$para7 = $tr.Paragraphs(7)
$start7 = $para7.Start

# Remove the "Async" and "/await " runs entirely (first 12 characters)
$prefix = $tr.Characters($start7, 12)
$prefix.Delete()

# Replace "при работе " with "Асинхронные " (keeps the same run / formatting)
$run1 = $tr.Characters($start7, 11)
$run1.Text = "Асинхронные "

# Skip over "CRUD " (unchanged) and split "операций" into "операци" + "и"
$afterRun1 = $start7 + ("Асинхронные ".Length)
$crudLen = "CRUD ".Length
$opStart = $afterRun1 + $crudLen

# Changing only the last character forces PowerPoint to split it into its own run
$lastChar = $tr.Characters($opStart + 7, 1)
$lastChar.Text = "и"

# --- Edit 3: paragraph 10 ---
# Before: "Авторизация с помощью " + "Microsoft Identity"
# After : "Microsoft " + "Identity"
$para10 = $tr.Paragraphs(10)
$start10 = $para10.Start

# Remove the "Авторизация с помощью " run entirely (first 22 characters)
$prefix10 = $tr.Characters($start10, 22)
$prefix10.Delete()

# Split "Microsoft Identity" into "Microsoft " + "Identity" (two runs)
$identityPart = $tr.Characters($start10 + 10, 8)
$identityPart.Text = "Identity"
